# Gestion des titres et des bullets en extractionTXT
#
# Paragraph 3 currently reads:
#   "1 – 2 " + bookmarkStart(_GoBack) + bookmarkEnd(_GoBack) + "- 3 : Principe du libre échange"
# (i.e. two runs split around the _GoBack bookmark).
#
# Target state: the heading becomes a single merged run, followed by a new
# block of explanatory paragraphs (with a couple of blank separator
# paragraphs and proofing marks around the mis-typed/foreign words), and the
# _GoBack bookmark is moved down to sit at the very end of the new content.

$d = $word.ActiveDocument
$wordMlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Heading paragraph: one run, no more embedded bookmark.
$heading = "<w:p $wordMlNs><w:r><w:t>1 – 2 - 3 : Principe du libre échange</w:t></w:r></w:p>"

$blank = "<w:p $wordMlNs/>"

# "On ne peut avoir une partie de la forme précèdente sans contenu après"
# with a spell-check proofErr wrapped around "précèdente".
$noPrecedingEmptyPart =
    "<w:p $wordMlNs>" +
        '<w:r><w:t xml:space="preserve">On ne peut avoir une partie de la forme </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>précèdente</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> sans contenu après</w:t></w:r>' +
    '</w:p>'

$mustIncrementLevel =
    "<w:p $wordMlNs><w:r><w:t>Si du contenu apparait, on est obligé d’incrémenter le niveau choisit après.</w:t></w:r></w:p>"

$noShapeLike =
    "<w:p $wordMlNs><w:r><w:t>On ne peut avoir qqch de la forme :</w:t></w:r></w:p>"

# "1-2 : theme 1" with proofErr wrapped around "theme".
$theme1 =
    "<w:p $wordMlNs>" +
        '<w:r><w:t xml:space="preserve">1-2 : </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>theme</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> 1</w:t></w:r>' +
    '</w:p>'

$slidesContenu1 =
    "<w:p $wordMlNs><w:r><w:t>Slides Contenu</w:t></w:r></w:p>"

# "1-2-1 : Theme 1.1" with proofErr wrapped around "Theme".
$theme11 =
    "<w:p $wordMlNs>" +
        '<w:r><w:t xml:space="preserve">1-2-1 : </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Theme</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> 1.1</w:t></w:r>' +
    '</w:p>'

# Last "Slides Contenu" paragraph now carries the relocated _GoBack bookmark.
$slidesContenu2 =
    "<w:p $wordMlNs>" +
        '<w:r><w:t>Slides Contenu</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

$newContentXml = $heading + $blank + $noPrecedingEmptyPart + $blank +
    $mustIncrementLevel + $noShapeLike + $theme1 + $slidesContenu1 +
    $theme11 + $slidesContenu2

# Replacing the whole paragraph-3 range (which spans both pre-bookmark and
# post-bookmark runs) with the new WordML fragment substitutes all of it in
# one shot, dropping the old mid-paragraph bookmark along the way.
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertXML($newContentXml)
